# Scene.xlsx fix: correct "Ini" -> "ini" casing in the FilePath column
# (F9:F14 on the DataNode sheet) and move the saved cell selection to G17,
# matching the author's final cursor position when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataNode")

$ws.Range("F9").Value  = "../resource/ini/Scene/1.xml"
$ws.Range("F10").Value = "../resource/ini/Scene/2.xml"
$ws.Range("F11").Value = "../resource/ini/Scene/3.xml"
$ws.Range("F12").Value = "../resource/ini/Scene/4.xml"
$ws.Range("F13").Value = "../resource/ini/Scene/5.xml"
$ws.Range("F14").Value = "../resource/ini/Scene/6.xml"

# Reflect the new active cell / selection saved in the workbook.
$null = $ws.Range("G17").Select()
